# Converting LunchMenu form html/java to xhtml/facelets/java
#
# Content changes applied to the "Blad1" lunch-menu sheet:
#   1. A new blank row is inserted above the current Monday row, pushing
#      the whole weekly menu down by one row (so it now lives in rows 2-6).
#      The new blank row keeps the row height the menu's first row used to
#      have.
#   2. A new "Lordag" (Saturday) menu entry is appended as a new last row:
#      Lördag / Grillad Ryggbiff / med pommes och hemmaslagen bea / 190
#      formatted the same way as the other weekday rows.
#   3. The active selection ends up on C13, matching where the editor left
#      the cursor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row above row 1; existing data shifts down to rows 2-6.
$ws.Rows.Item(1).Insert()
# Restore the original header row's height on the new blank row.
$ws.Rows.Item(1).RowHeight = 22.8

# 2. Add the new Saturday menu row at the bottom (now row 7).
$ws.Range("A7").Value = "Lördag"
$ws.Range("B7").Value = "Grillad Ryggbiff"
$ws.Range("C7").Value = "med pommes och hemmaslagen bea"
$ws.Range("D7").Value = 190

# Match the formatting (vertical-centered, wrapped text) used by the other
# menu rows by copying the format from the row right above the new one.
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 3. Leave the selection where the editor left it.
$ws.Range("C13").Select() | Out-Null
